# Hide slides 5, 7, 8, 9, 10, 11 (1-based slide index in the slide show
# order) in the presentation, matching the author's "hide slides in
# moreselectors ppt" commit.

$p = $ppt.ActivePresentation

$hiddenIndexes = @(5, 7, 8, 9, 10, 11)

foreach ($idx in $hiddenIndexes) {
    $s = $p.Slides.Item($idx)
    $s.SlideShowTransition.Hidden = $true
}
